$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*checkboxs*") {
        $paraStart = $p.Range.Start
        $paraEnd = $p.Range.End

        # Locate the end of "formulário " so the preceding run is left untouched.
        $anchor = $d.Range($paraStart, $paraEnd)
        $anchor.Find.Execute("formulário ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

        $target = $d.Range($anchor.End, $paraEnd)
        $target.Find.ClearFormatting()
        $target.Find.Replacement.ClearFormatting()
        $target.Find.Execute(
            "com campo inserção de palavra chave para pesquisa e checkboxs para filtrá-las.",
            $true, $false, $false, $false, $false,
            $true, 1, $false,
            "com campo inserção de palavra chave para pesquisa.",
            2
        )
    }
}
